$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Il1rapl1"
$ws.Cells.Item(2, 3).Value = "Ptprf"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.03825666666666667
$ws.Cells.Item(2, 8).Value = 0.11477
$ws.Cells.Item(2, 9).Value = 0.4331352014340976
$ws.Cells.Item(2, 10).Value = 0.4331352014340976
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.242595
$ws.Cells.Item(2, 14).Value = 0.727785
$ws.Cells.Item(2, 15).Value = 0.03230700759563258
$ws.Cells.Item(2, 16).Value = 0.03230700759563257
$ws.Cells.Item(2, 17).Value = 0.009280876050000001
$ws.Cells.Item(2, 18).Value = 0.08352788445000001
$ws.Cells.Item(2, 19).Value = 0.01399330224266724
$ws.Cells.Item(2, 20).Value = 0.01399330224266723

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Il1rapl1"
$ws.Cells.Item(3, 3).Value = "Ptprf"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.03825666666666667
$ws.Cells.Item(3, 8).Value = 0.11477
$ws.Cells.Item(3, 9).Value = 0.4331352014340976
$ws.Cells.Item(3, 10).Value = 0.4331352014340976
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.674351333333334
$ws.Cells.Item(3, 14).Value = 11.023054
$ws.Cells.Item(3, 15).Value = 0.4893229309549773
$ws.Cells.Item(3, 16).Value = 0.4893229309549771
$ws.Cells.Item(3, 17).Value = 0.1405684341755556
$ws.Cells.Item(3, 18).Value = 1.26511590758
$ws.Cells.Item(3, 19).Value = 0.2119429862655071
$ws.Cells.Item(3, 20).Value = 0.211942986265507

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Il1rapl1"
$ws.Cells.Item(4, 3).Value = "Ptprf"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.03825666666666667
$ws.Cells.Item(4, 8).Value = 0.11477
$ws.Cells.Item(4, 9).Value = 0.4331352014340976
$ws.Cells.Item(4, 10).Value = 0.4331352014340976
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.574634666666667
$ws.Cells.Item(4, 14).Value = 10.723904
$ws.Cells.Item(4, 15).Value = 0.4760434029044768
$ws.Cells.Item(4, 16).Value = 0.4760434029044767
$ws.Cells.Item(4, 17).Value = 0.1367536068977778
$ws.Cells.Item(4, 18).Value = 1.23078246208
$ws.Cells.Item(4, 19).Value = 0.2061911552084038
$ws.Cells.Item(4, 20).Value = 0.2061911552084038

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Il1rapl1"
$ws.Cells.Item(5, 3).Value = "Ptprf"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.03825666666666667
$ws.Cells.Item(5, 8).Value = 0.11477
$ws.Cells.Item(5, 9).Value = 0.4331352014340976
$ws.Cells.Item(5, 10).Value = 0.4331352014340976
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.017471
$ws.Cells.Item(5, 14).Value = 0.052413
$ws.Cells.Item(5, 15).Value = 0.002326658544913526
$ws.Cells.Item(5, 16).Value = 0.002326658544913525
$ws.Cells.Item(5, 17).Value = 0.0006683822233333334
$ws.Cells.Item(5, 18).Value = 0.00601544001
$ws.Cells.Item(5, 19).Value = 0.001007757717519485
$ws.Cells.Item(5, 20).Value = 0.001007757717519484

# Row 6
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Il1rapl1"
$ws.Cells.Item(6, 3).Value = "Ptprf"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.05006833333333333
$ws.Cells.Item(6, 8).Value = 0.150205
$ws.Cells.Item(6, 9).Value = 0.5668647985659024
$ws.Cells.Item(6, 10).Value = 0.5668647985659024
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.242595
$ws.Cells.Item(6, 14).Value = 0.727785
$ws.Cells.Item(6, 15).Value = 0.03230700759563258
$ws.Cells.Item(6, 16).Value = 0.03230700759563257
$ws.Cells.Item(6, 17).Value = 0.012146327325
$ws.Cells.Item(6, 18).Value = 0.109316945925
$ws.Cells.Item(6, 19).Value = 0.01831370535296534
$ws.Cells.Item(6, 20).Value = 0.01831370535296534

# Row 7
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Il1rapl1"
$ws.Cells.Item(7, 3).Value = "Ptprf"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.05006833333333333
$ws.Cells.Item(7, 8).Value = 0.150205
$ws.Cells.Item(7, 9).Value = 0.5668647985659024
$ws.Cells.Item(7, 10).Value = 0.5668647985659024
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.674351333333334
$ws.Cells.Item(7, 14).Value = 11.023054
$ws.Cells.Item(7, 15).Value = 0.4893229309549773
$ws.Cells.Item(7, 16).Value = 0.4893229309549771
$ws.Cells.Item(7, 17).Value = 0.1839686473411111
$ws.Cells.Item(7, 18).Value = 1.65571782607
$ws.Cells.Item(7, 19).Value = 0.2773799446894702
$ws.Cells.Item(7, 20).Value = 0.2773799446894701

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Il1rapl1"
$ws.Cells.Item(8, 3).Value = "Ptprf"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.05006833333333333
$ws.Cells.Item(8, 8).Value = 0.150205
$ws.Cells.Item(8, 9).Value = 0.5668647985659024
$ws.Cells.Item(8, 10).Value = 0.5668647985659024
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 3.574634666666667
$ws.Cells.Item(8, 14).Value = 10.723904
$ws.Cells.Item(8, 15).Value = 0.4760434029044768
$ws.Cells.Item(8, 16).Value = 0.4760434029044767
$ws.Cells.Item(8, 17).Value = 0.1789760000355555
$ws.Cells.Item(8, 18).Value = 1.61078400032
$ws.Cells.Item(8, 19).Value = 0.269852247696073
$ws.Cells.Item(8, 20).Value = 0.2698522476960729

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Il1rapl1"
$ws.Cells.Item(9, 3).Value = "Ptprf"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.05006833333333333
$ws.Cells.Item(9, 8).Value = 0.150205
$ws.Cells.Item(9, 9).Value = 0.5668647985659024
$ws.Cells.Item(9, 10).Value = 0.5668647985659024
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.017471
$ws.Cells.Item(9, 14).Value = 0.052413
$ws.Cells.Item(9, 15).Value = 0.002326658544913526
$ws.Cells.Item(9, 16).Value = 0.002326658544913525
$ws.Cells.Item(9, 17).Value = 0.0008747438516666665
$ws.Cells.Item(9, 18).Value = 0.007872694664999999
$ws.Cells.Item(9, 19).Value = 0.001318900827394042
$ws.Cells.Item(9, 20).Value = 0.001318900827394041
